# Update the three subscriber email addresses, stripping their old
# "mailto:" hyperlinks/rich-text underline formatting in the process.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all hyperlinks from the sheet (also drops the underlined,
# colored rich-text run that previously rendered each email).
$ws.Hyperlinks.Delete()

# Replace the old addresses with the new ones as plain text.
$ws.Range("C2").Value = "aegxhpzio@yomail.info"
$ws.Range("C3").Value = "xlvupdwec@firste.ml"
$ws.Range("C4").Value = "aegxhqbpc@yomail.info"

# Leave the final selection on C9, matching the saved workbook state.
$ws.Range("C9").Select()
